$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.474.80'
$ws.Range('E2').Value = '  +0.11%  '

$ws.Range('D3').Value = '1.824.89'
$ws.Range('E3').Value = '  -0.15%  '

$ws.Range('D4').Value = '''1.005'
$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').Value = '''316.23'
$ws.Range('E5').Value = '  +0.48%  '

$ws.Range('E6').Value = '  +0.23%  '

$ws.Range('D7').Value = '''0.5170'
$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').Value = '''0.3853'
$ws.Range('E8').Value = '  -1.26%  '

$ws.Range('D9').Value = '''0.08276'
$ws.Range('E9').Value = '  +8.24%  '

$ws.Range('D10').Value = '''1.123'
$ws.Range('E10').Value = '  +1.26%  '

$ws.Range('E11').Value = '  +0.02%  '

$ws.Range('D12').Value = '''6.372'
$ws.Range('E12').Value = '  +1.52%  '

$ws.Range('D13').Value = '''21.07'
$ws.Range('E13').Value = '  +0.26%  '

$ws.Range('E14').Value = '  +0.27%  '

$ws.Range('D15').Value = '''7.470'
$ws.Range('E15').Value = '  -0.88%  '

$ws.Range('D16').Value = '1.818.53'
$ws.Range('E16').Value = '  -0.27%  '

$ws.Range('D17').Value = '''94.06'
$ws.Range('E17').Value = '  +0.78%  '

$ws.Range('E18').Value = '  +3.90%  '

$ws.Range('D19').Value = '''0.06630'
$ws.Range('E19').Value = '  -0.54%  '

$ws.Range('E20').Value = '  +0.75%  '

$ws.Range('D21').Value = '''1.004'
$ws.Range('E21').Value = '  +0.29%  '

$ws.Range('D22').Value = '''6.045'
$ws.Range('E22').Value = '  -2.00%  '

$ws.Range('D23').Value = '28.507.21'
$ws.Range('E23').Value = '  +0.13%  '

$ws.Range('D24').Value = '''11.55'
$ws.Range('E24').Value = '  +3.46%  '

$ws.Range('D25').Value = '''2.244'
$ws.Range('E25').Value = '  -0.60%  '

$ws.Range('D26').Value = '''21.06'
$ws.Range('E26').Value = '  +2.18%  '

$ws.Range('D27').Value = '''159.34'
$ws.Range('E27').Value = '  +1.44%  '

$ws.Range('D28').Value = '2.032.35'
$ws.Range('E28').Value = '  -0.11%  '

$ws.Range('D29').Value = '''2.403'
$ws.Range('E29').Value = '  +0.26%  '

$ws.Range('D30').Value = '''125.68'
$ws.Range('E30').Value = '  +0.73%  '

$ws.Range('E31').Value = '  +2.15%  '

$ws.Range('D32').Value = '''1.096'
$ws.Range('E32').Value = '  -2.00%  '

$ws.Range('E33').Value = '  +1.37%  '

$ws.Range('D34').Value = '''0.07525'
$ws.Range('E34').Value = '  +7.47%  '

$ws.Range('E35').Value = '  +0.44%  '

$ws.Range('D36').Value = '''12.25'
$ws.Range('E36').Value = '  +9.27%  '

$ws.Range('D37').Value = '''0.2227'
$ws.Range('E37').Value = '  +0.21%  '

$ws.Range('D38').Value = '''0.02363'
$ws.Range('E38').Value = '  +1.85%  '

$ws.Range('D39').Value = '''5.229'
$ws.Range('E39').Value = '  +2.04%  '

$ws.Range('D40').Value = '''8.772'
$ws.Range('E40').Value = '  -2.37%  '

$ws.Range('D41').Value = '''0.6387'
$ws.Range('E41').Value = '  +1.71%  '

$ws.Range('D42').Value = '''1.188'
$ws.Range('E42').Value = '  +0.06%  '

$ws.Range('E43').Value = '  -0.19%  '

$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').Value = '''0.6205'
$ws.Range('E44').Value = '  +5.24%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''13.67'
$ws.Range('E45').Value = '  +1.86%  '

$ws.Range('D46').Value = '''3.804'
$ws.Range('E46').Value = '  +2.42%  '

$ws.Range('D47').Value = '''127.73'
$ws.Range('E47').Value = '  +2.89%  '

$ws.Range('D48').Value = '''2.018'
$ws.Range('E48').Value = '  +2.11%  '

$ws.Range('E49').Value = '  +0.48%  '

$ws.Range('D50').Value = '''0.06950'
$ws.Range('E50').Value = '  +0.39%  '

$ws.Range('E51').Value = '  +1.58%  '
